# Remove the trailing "Ver no Jupiter ..." / "(c) 2020 ..." footer block
# that the Jekyll site build appended, along with the blank paragraph that
# used to separate the two, leaving a single blank paragraph before the
# page-break paragraph -- matching the upstream commit that rebuilt the
# site without this footer.

$d = $word.ActiveDocument

$target1 = "Ver no Jupiter Salvar em pdf Salvar em docx"
$target2 = "$([char]0x00A9) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13)
    if ($text -eq $target1) {
        $startPara = $p
    }
    if ($text -eq $target2) {
        $endPara = $p
        break
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    # Extend the end of the deletion through the following blank paragraph
    # (there were two blank paragraphs surrounding this footer block; only
    # one should remain afterwards), i.e. delete from the start of the
    # "Ver no Jupiter..." paragraph through the end of the blank paragraph
    # that immediately follows the copyright paragraph.
    $trailingBlank = $endPara.Next()
    $deleteEnd = $endPara.Range.End
    if ($trailingBlank -ne $null -and $trailingBlank.Range.Text.TrimEnd([char]13) -eq "") {
        $deleteEnd = $trailingBlank.Range.End
    }

    $r = $d.Range($startPara.Range.Start, $deleteEnd)
    $r.Delete()
}

Write-Output "Paragraphs remaining: $($d.Paragraphs.Count)"
